# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Anima_Profits leve-profit columns (H:N)
# across several worksheets, per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("M62").Value = $null

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("M65").Value = $null

$ws.Range("H98").Value = 1018.1667
$ws.Range("I98").Value = 690.8889
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 690.8889
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 807.1111
$ws.Range("N98").Value = -4996

$ws.Range("H106").Value = 17244390
$ws.Range("I106").Value = 21742100
$ws.Range("J106").Value = 3166.6667
$ws.Range("K106").Value = 21742100
$ws.Range("L106").Value = 3166.6667
$ws.Range("M106").Value = -21741469
$ws.Range("N106").Value = -4428.6667

$ws.Range("H122").Value = 1018.1667
$ws.Range("I122").Value = 690.8889
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2072.6667
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 377.3332999999998
$ws.Range("N122").Value = -10900

$ws.Range("H137").Value = 1149.8853
$ws.Range("I137").Value = 973.6889
$ws.Range("K137").Value = 2921.0667
$ws.Range("M137").Value = -371.0666999999999

$ws.Range("H138").Value = 1146.18
$ws.Range("I138").Value = 494.01923
$ws.Range("J138").Value = 1852.6875
$ws.Range("K138").Value = 1482.05769
$ws.Range("L138").Value = 5558.0625
$ws.Range("M138").Value = 3657.94231
$ws.Range("N138").Value = -15838.0625

$ws.Range("H141").Value = 2602.8965
$ws.Range("I141").Value = 881.7727
$ws.Range("J141").Value = 8012.143
$ws.Range("K141").Value = 2645.3181
$ws.Range("L141").Value = 24036.429
$ws.Range("M141").Value = 2534.6819
$ws.Range("N141").Value = -34396.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 5000
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5350

$ws.Range("H32").Value = 913078.0600000001
$ws.Range("I32").Value = 1047077.6
$ws.Range("J32").Value = 19747.445
$ws.Range("K32").Value = 1047077.6
$ws.Range("L32").Value = 19747.445
$ws.Range("M32").Value = -1046790.6
$ws.Range("N32").Value = -20321.445

$ws.Range("H61").Value = 1964.5714
$ws.Range("I61").Value = 1655.0555
$ws.Range("K61").Value = 1655.0555
$ws.Range("M61").Value = -1443.0555

$ws.Range("H74").Value = 902.4262
$ws.Range("I74").Value = 681
$ws.Range("J74").Value = 1294.9546
$ws.Range("K74").Value = 681
$ws.Range("L74").Value = 1294.9546
$ws.Range("M74").Value = 193
$ws.Range("N74").Value = -3042.9546

$ws.Range("H76").Value = 32000
$ws.Range("J76").Value = 32000
$ws.Range("L76").Value = 32000
$ws.Range("N76").Value = -32676

$ws.Range("H77").Value = 902.4262
$ws.Range("I77").Value = 681
$ws.Range("J77").Value = 1294.9546
$ws.Range("K77").Value = 3405
$ws.Range("L77").Value = 6474.773
$ws.Range("M77").Value = 963
$ws.Range("N77").Value = -15210.773

$ws.Range("H79").Value = 32000
$ws.Range("J79").Value = 32000
$ws.Range("L79").Value = 32000
$ws.Range("N79").Value = -34340

$ws.Range("H136").Value = 1964.5714
$ws.Range("I136").Value = 1655.0555
$ws.Range("K136").Value = 4965.166499999999
$ws.Range("M136").Value = -2415.166499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2085.9473
$ws.Range("I80").Value = 2202.4443
$ws.Range("J80").Value = 1981.1
$ws.Range("K80").Value = 2202.4443
$ws.Range("L80").Value = 1981.1
$ws.Range("M80").Value = -1204.4443
$ws.Range("N80").Value = -3977.1

$ws.Range("H82").Value = 26734.385
$ws.Range("I82").Value = 6714.25
$ws.Range("J82").Value = 35632.223
$ws.Range("K82").Value = 6714.25
$ws.Range("L82").Value = 35632.223
$ws.Range("M82").Value = -6331.25
$ws.Range("N82").Value = -36398.223

$ws.Range("H83").Value = 2085.9473
$ws.Range("I83").Value = 2202.4443
$ws.Range("J83").Value = 1981.1
$ws.Range("K83").Value = 11012.2215
$ws.Range("L83").Value = 9905.5
$ws.Range("M83").Value = -6020.2215
$ws.Range("N83").Value = -19889.5

$ws.Range("H85").Value = 26734.385
$ws.Range("I85").Value = 6714.25
$ws.Range("J85").Value = 35632.223
$ws.Range("K85").Value = 6714.25
$ws.Range("L85").Value = 35632.223
$ws.Range("M85").Value = -5388.25
$ws.Range("N85").Value = -38284.223

$ws.Range("H134").Value = 2461.4792
$ws.Range("I134").Value = 2245.6553
$ws.Range("J134").Value = 2790.8948
$ws.Range("K134").Value = 6736.965899999999
$ws.Range("L134").Value = 8372.6844
$ws.Range("M134").Value = -4201.965899999999
$ws.Range("N134").Value = -13442.6844

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 990
$ws.Range("I16").Value = 980
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 980
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -693
$ws.Range("N16").Value = -1574

$ws.Range("H31").Value = 4666.362
$ws.Range("I31").Value = 1347.3478
$ws.Range("J31").Value = 7847.0835
$ws.Range("K31").Value = 1347.3478
$ws.Range("L31").Value = 7847.0835
$ws.Range("M31").Value = -1052.3478
$ws.Range("N31").Value = -8437.083500000001

$ws.Range("H34").Value = 4666.362
$ws.Range("I34").Value = 1347.3478
$ws.Range("J34").Value = 7847.0835
$ws.Range("K34").Value = 1347.3478
$ws.Range("L34").Value = 7847.0835
$ws.Range("M34").Value = -1145.3478
$ws.Range("N34").Value = -8251.083500000001

$ws.Range("H58").Value = 938.1731
$ws.Range("I58").Value = 637.8788
$ws.Range("J58").Value = 1459.7368
$ws.Range("K58").Value = 637.8788
$ws.Range("L58").Value = 1459.7368
$ws.Range("M58").Value = -434.8788
$ws.Range("N58").Value = -1865.7368

$ws.Range("H68").Value = 23866.111
$ws.Range("J68").Value = 23866.111
$ws.Range("L68").Value = 23866.111
$ws.Range("N68").Value = -25364.111

$ws.Range("H71").Value = 23866.111
$ws.Range("J71").Value = 23866.111
$ws.Range("L71").Value = 71598.333
$ws.Range("N71").Value = -79086.333

$ws.Range("H99").Value = 1784.8
$ws.Range("I99").Value = 1731
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1731
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -233
$ws.Range("N99").Value = -4996

$ws.Range("H113").Value = 990
$ws.Range("I113").Value = 980
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 980
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1190
$ws.Range("N113").Value = -5340

$ws.Range("H122").Value = 1497.3334
$ws.Range("I122").Value = 1596
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 4788
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -2338
$ws.Range("N122").Value = -8800

$ws.Range("H126").Value = 1784.8
$ws.Range("I126").Value = 1731
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5193
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2723
$ws.Range("N126").Value = -10940

$ws.Range("H132").Value = 4168096.2
$ws.Range("I132").Value = 1392.9166
$ws.Range("J132").Value = 10418151
$ws.Range("K132").Value = 4178.7498
$ws.Range("L132").Value = 31254453
$ws.Range("M132").Value = -1648.7498
$ws.Range("N132").Value = -31259513

$ws.Range("H134").Value = 2509.077
$ws.Range("I134").Value = 2743.0188
$ws.Range("J134").Value = 1475.8334
$ws.Range("K134").Value = 8229.056399999999
$ws.Range("L134").Value = 4427.5002
$ws.Range("M134").Value = -5694.056399999999
$ws.Range("N134").Value = -9497.5002

$ws.Range("H136").Value = 938.1731
$ws.Range("I136").Value = 637.8788
$ws.Range("J136").Value = 1459.7368
$ws.Range("K136").Value = 1913.6364
$ws.Range("L136").Value = 4379.2104
$ws.Range("M136").Value = 636.3636000000001
$ws.Range("N136").Value = -9479.2104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1666.6666
$ws.Range("J75").Value = 1666.6666
$ws.Range("L75").Value = 4999.9998
$ws.Range("N75").Value = -6995.9998

$ws.Range("H78").Value = 1666.6666
$ws.Range("J78").Value = 1666.6666
$ws.Range("L78").Value = 14999.9994
$ws.Range("N78").Value = -24983.9994

$ws.Range("H131").Value = 2716.1572
$ws.Range("J131").Value = 2884.9692
$ws.Range("L131").Value = 8654.9076
$ws.Range("N131").Value = -18734.9076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 500
$ws.Range("K6").Value = 500
$ws.Range("M6").Value = -387

$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("K16").Value = 500
$ws.Range("M16").Value = -250

$ws.Range("H122").Value = 2750
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 169067.33
$ws.Range("I40").Value = 252101
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 252101
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -251965
$ws.Range("N40").Value = -3272

$ws.Range("H132").Value = 2464.6416
$ws.Range("I132").Value = 2348.111
$ws.Range("J132").Value = 2711.4119
$ws.Range("K132").Value = 7044.333
$ws.Range("L132").Value = 8134.2357
$ws.Range("M132").Value = -4514.333
$ws.Range("N132").Value = -13194.2357

$ws.Range("H136").Value = 4903568.5
$ws.Range("I136").Value = 1585.5385
$ws.Range("J136").Value = 20835012
$ws.Range("K136").Value = 4756.6155
$ws.Range("L136").Value = 62505036
$ws.Range("M136").Value = -2206.6155
$ws.Range("N136").Value = -62510136

$ws.Range("H17").Value = 70005
$ws.Range("J17").Value = 70005
$ws.Range("L17").Value = 70005
$ws.Range("N17").Value = -70349

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2689336
$ws.Range("I132").Value = 1261.0488
$ws.Range("J132").Value = 7937482.5
$ws.Range("K132").Value = 3783.1464
$ws.Range("L132").Value = 23812447.5
$ws.Range("M132").Value = -1253.1464
$ws.Range("N132").Value = -23817507.5
